$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range("M4").Value = "Golden Eagle"
$ws.Range("F8").Value = "Okapi"
$ws.Range("E12").Value = "Striped dolphin"
$ws.Range("D14").Value = "Striped dolphin"
$ws.Range("G16").Value = "Okapi"
$ws.Range("N18").Value = "Goanna"
$ws.Range("M20").Value = "Homo habilis"
$ws.Range("N22").Value = "Homo habilis"
$ws.Range("F24").Value = "Kudu"
$ws.Range("D26").Value = "Striped Rabbit"
$ws.Range("E28").Value = "Kudu"
$ws.Range("D30").Value = "Kudu"
$ws.Range("E36").Value = "Sea Otter"
$ws.Range("F40").Value = "Sea Otter"
$ws.Range("L40").Value = "Pacific Spiny Lumpsucker"
$ws.Range("N42").Value = "Pacific Spiny Lumpsucker"
$ws.Range("E44").Value = "Siberian Chipmunk"
$ws.Range("M44").Value = "Pacific Spiny Lumpsucker"
$ws.Range("D46").Value = "Siberian Chipmunk"
$ws.Range("N46").Value = "Darwin's Frogs"
$ws.Range("D50").Value = "Itjaritjari"
$ws.Range("N50").Value = "Bat-Eared Fox"
